$wb = $excel.ActiveWorkbook

# --- Character sheet: add the "Skill" columns data (P:T) for the first two rows ---
$wsChar = $wb.Worksheets.Item("Character")

# Row 2 (ID 1000): SkillID1..SkillID4 -> 9001000, patternID -> 8001000
$wsChar.Range("P2:S2").Value = 9001000
$wsChar.Range("T2").Value = 8001000

# Row 3 (ID 2000): new Skill/pattern values, matching the formatting used in row 2
$wsChar.Range("P2:T2").Copy()
$wsChar.Range("P3:T3").PasteSpecial(-4122)
$wsChar.Range("P3:S3").Value = 9002000
$wsChar.Range("T3").Value = 8001000

# --- Selection / active-tab bookkeeping, matching the saved workbook state ---
$wsEnemy = $wb.Worksheets.Item("Enemy")
[void]$wsEnemy.Range("F1:M2").Select()

[void]$wsChar.Activate()
[void]$wsChar.Range("T4").Select()
